$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.855.91"
$ws.Range("E2").Value = "  -1.27%  "
$ws.Range("D3").Value = "2.327.49"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.18"
$ws.Range("E5").Value = "  -1.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.26"
$ws.Range("E6").Value = "  -2.96%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.509"
$ws.Range("E7").Value = "  -4.48%  "
$ws.Range("E9").Value = "  -4.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.50"
$ws.Range("E10").Value = "  -5.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.16"
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("E12").Value = "  -2.27%  "
$ws.Range("E13").Value = "  +0.70%  "
$ws.Range("E14").Value = "  -3.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.79"
$ws.Range("E15").Value = "  +4.33%  "
$ws.Range("D16").Value = "2.330.08"
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("E17").Value = "  +1.80%  "
$ws.Range("D18").Value = "42.778.09"
$ws.Range("E18").Value = "  -1.19%  "
$ws.Range("E20").Value = "  -2.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.60"
$ws.Range("E21").Value = "  -5.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.13"
$ws.Range("E22").Value = "  +1.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.90"
$ws.Range("E23").Value = "  -2.99%  "
$ws.Range("E24").Value = "  -2.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.54"
$ws.Range("E25").Value = "  -2.85%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.33"
$ws.Range("E27").Value = "  +2.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.96"
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("E29").Value = "  -5.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.73"
$ws.Range("E30").Value = "  -6.11%  "
$ws.Range("E31").Value = "  -4.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "161.13"
$ws.Range("E32").Value = "  -4.37%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("E34").Value = "  -4.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.60"
$ws.Range("E35").Value = "  +2.71%  "
$ws.Range("E36").Value = "  -3.40%  "
$ws.Range("B37").Value = "Celestia"
$ws.Range("C37").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.12"
$ws.Range("E37").Value = "  -5.49%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0717"
$ws.Range("E38").Value = "  -3.55%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.88"
$ws.Range("E39").Value = "  -5.48%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.83"
$ws.Range("E40").Value = "  -2.59%  "
$ws.Range("E42").Value = "  -3.23%  "
$ws.Range("E43").Value = "  +0.42%  "
$ws.Range("D44").Value = "2.009.21"
$ws.Range("E44").Value = "  +1.24%  "
$ws.Range("E45").Value = "  -4.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "18.53"
$ws.Range("E46").Value = "  -2.64%  "
$ws.Range("E47").Value = "  +1.70%  "
$ws.Range("E48").Value = "  -4.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.49"
$ws.Range("E49").Value = "  -0.93%  "
$ws.Range("E50").Value = "  -1.63%  "
$ws.Range("D51").Value = "2.556.40"
$ws.Range("E51").Value = "  +0.93%  "
